# Auto update stock data
# Updates the "Date_1" (column A) and "EBITDA" (column B) values for the
# most-recent-date rows of each company block, advancing the snapshot date
# from 2025/12/05 to 2025/12/06 and refreshing the EBITDA figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> [Date, EBITDA] (EBITDA left blank/null when unchanged)
$updates = @(
    @{ Row = 2;  Date = "2025/12/06"; Ebitda = "5.36" },
    @{ Row = 8;  Date = "2025/12/06"; Ebitda = "7.88" },
    @{ Row = 14; Date = "2025/12/06"; Ebitda = $null },
    @{ Row = 20; Date = "2025/12/06"; Ebitda = "12.63" },
    @{ Row = 26; Date = "2025/12/06"; Ebitda = "10.50" },
    @{ Row = 32; Date = "2025/12/06"; Ebitda = "26.61" },
    @{ Row = 38; Date = "2025/12/06"; Ebitda = $null },
    @{ Row = 44; Date = "2025/12/06"; Ebitda = "11.27" },
    @{ Row = 50; Date = "2025/12/06"; Ebitda = "12.12" },
    @{ Row = 56; Date = "2025/12/06"; Ebitda = "33.94" },
    @{ Row = 62; Date = "2025/12/06"; Ebitda = "11.93" },
    @{ Row = 68; Date = "2025/12/06"; Ebitda = "13.00" },
    @{ Row = 74; Date = "2025/12/06"; Ebitda = "16.67" }
)

foreach ($u in $updates) {
    $rowNum = $u.Row

    # Column A holds text dates like "2025/12/05" - force text format so
    # Excel does not reinterpret the string as a date serial number.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $u.Date

    if ($null -ne $u.Ebitda) {
        # Column B holds text numbers like "5.39" - force text format so
        # Excel keeps the exact string instead of converting to a Double.
        $cellB = $ws.Cells.Item($rowNum, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $u.Ebitda
    }
}
